$wb = $excel.ActiveWorkbook

# --- "Budget Out" sheet: row 9 amount + description text ---
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$wsBudgetOut.Range("C9").Value = 89.02

# --- "TestRecord" sheet: row 10 date + amount + note text ---
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$wsTestRecord.Range("A10").Value = 43262
$wsTestRecord.Range("B10").Value = 119.94

# --- "Expected Out" sheet: rows 9, 11 amounts (B1 SUM formula recalculates automatically) ---
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1346.96
$wsExpectedOut.Range("B11").Value = 426.82

$excel.Calculate()
